# Tnfsf13-Tnfrsf1a.xlsx was regenerated with updated TPM input values.
# The ligand (Tnfsf13) average/total expression for the "ECs" sending
# cluster and the receptor (Tnfrsf1a) average/total expression for the
# "ECs" target cluster changed, which ripples into every derived
# specificity / edge-weight column (I, J, O, P, Q, R, S, T) for every
# sending x target cluster combination (rows 2-10).
#
# Values below are taken from the recomputed NATMI output (new TPM run);
# only the cells that actually differ from the original values are
# written here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New per-row values, keyed by row number -> column letter -> value.
$updates = [ordered]@{
    2  = @{
        G = 0.5347833333333334; H = 1.60435
        I = 0.196822066153855;  J = 0.196822066153855
        M = 19.98610666666667;  N = 59.95832
        O = 0.2969043109767812; P = 0.2969043109767812
        Q = 10.68823674355556;  R = 96.19413069200002
        S = 0.05843731993643678; T = 0.05843731993643677
    }
    3  = @{
        G = 0.5347833333333334; H = 1.60435
        I = 0.196822066153855;  J = 0.196822066153855
        O = 0.4664722083712238; P = 0.4664722083712239
        Q = 16.79249917577222;  R = 151.13249258195
        S = 0.09181202385497585; T = 0.09181202385497586
    }
    4  = @{
        G = 0.5347833333333334; H = 1.60435
        I = 0.196822066153855;  J = 0.196822066153855
        O = 0.236623480651995;  P = 0.236623480651995
        Q = 8.518191507466668;  R = 76.66372356720001
        S = 0.04657272236244239; T = 0.04657272236244239
    }
    5  = @{
        I = 0.1891972429821067; J = 0.1891972429821067
        M = 19.98610666666667;  N = 59.95832
        O = 0.2969043109767812; P = 0.2969043109767812
        Q = 10.27417790970667;  R = 92.46760118736
        S = 0.05617347706630905; T = 0.05617347706630904
    }
    6  = @{
        I = 0.1891972429821067; J = 0.1891972429821067
        O = 0.4664722083712238; P = 0.4664722083712239
        S = 0.08825525575161033; T = 0.08825525575161035
    }
    7  = @{
        I = 0.1891972429821067; J = 0.1891972429821067
        O = 0.236623480651995;  P = 0.236623480651995
        S = 0.04476851016418731; T = 0.04476851016418732
    }
    8  = @{
        I = 0.6139806908640383; J = 0.6139806908640382
        M = 19.98610666666667;  N = 59.95832
        O = 0.2969043109767812; P = 0.2969043109767812
        Q = 33.34164257170667;  R = 300.07478314536
        S = 0.1822935139740354; T = 0.1822935139740354
    }
    9  = @{
        I = 0.6139806908640383; J = 0.6139806908640382
        O = 0.4664722083712238; P = 0.4664722083712239
        S = 0.2864049287646376; T = 0.2864049287646376
    }
    10 = @{
        I = 0.6139806908640383; J = 0.6139806908640382
        O = 0.236623480651995;  P = 0.236623480651995
        S = 0.1452822481253652; T = 0.1452822481253653
    }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
